$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to stay plain text while we
# overwrite them, then restore the original (Normal) cell style so the
# on-disk formatting is unchanged -- only the text content differs, exactly
# like the source diff (t="inlineStr" cells, no "s" style attribute).
$colD = $ws.Range("D2:D51")
$colE = $ws.Range("E2:E51")
$colD.NumberFormat = "@"
$colE.NumberFormat = "@"

$ws.Range("D2").Value = "44.819.71"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "2.248.58"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "306.18"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.573"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -1.46%  "
$ws.Range("D10").Value = "35.34"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D12").Value = "7.23"
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("D14").Value = "2.590.99"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "0.842"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "2.241.58"
$ws.Range("E16").Value = "  -6.01%  "
$ws.Range("D17").Value = "13.59"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "44.596.17"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "11.95"
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").Value = "65.52"
$ws.Range("D23").Value = "240.08"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").Value = "9.88"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "37.54"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").Value = "6.03"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "150.84"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "0.0800"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "3.09"
$ws.Range("E35").Value = "  -6.58%  "
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "1.85"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D39").Value = "15.10"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "3.82"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.839.96"
$ws.Range("E44").Value = "  +5.79%  "
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  +13.28%  "
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "79.96"
$ws.Range("E47").Value = "  -5.55%  "
$ws.Range("D48").Value = "99.18"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").Value = "69.66"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "54.81"
$ws.Range("E51").Value = "  +0.79%  "

# Restore original style (remove the temporary text NumberFormat again).
$colD.Style = "Normal"
$colE.Style = "Normal"
